# Applies the 2025 cost-sheet update:
#  - Raises several annual base costs (Luz, Impuestos, Monotributo, seguro de
#    vida, limpieza, combustible)
#  - Converts row 30 (seguro de vida y accidentes) from a per-employee
#    formula to a plain annual value, bumps the per-unit rate in B30, and
#    documents the employee count in a new column F note
#  - Row 6's monthly formula is rewritten as an explicit =D6/12 (breaking
#    it out of the E3:E7 shared-formula group)
#  - Moves the sheet's scroll position / selection to reflect where the
#    user was working
#  - Widens a new column F to hold the "(8 empleados)" note

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 1: LOCAL -------------------------------------------------
$ws.Range("D6").Value = 3600000
$ws.Range("E6").Formula = "=D6/12"

$ws.Range("D7").Value = 360000

# --- Section 5: GASTOS PROFESIONALES ----------------------------------
$ws.Range("D29").Value = 1355000

$ws.Range("B30").Value = 2700
$ws.Range("D30").Value = 259200
$ws.Range("F30").Value = "(8 empleados)"

# --- Section 6: Mano de obra indirecta --------------------------------
$ws.Range("D34").Value = 1370000

# --- Section 7: Rodado --------------------------------------------------
$ws.Range("D36").Value = 2160000

# --- Column widths / view state ---------------------------------------
$ws.Columns.Item(6).ColumnWidth = 13.109375

$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("A35:E35").Select
